# Hortaliza, Vega Modelo de Temuco - Pimiento: insert 4 new daily price rows
# (week of 44610) above the existing row 1135 block, pushing the remaining
# historical rows down by 4 (old 1135-1212 -> new 1139-1216).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 1135, shifting rows 1135:1212 down to 1139:1216.
$ws.Rows("1135:1138").Insert()

# Common columns shared by every row of this market/category block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112002
$categoria = "Pimiento"
$clasif    = "Hortaliza"

$newRows = @(
    @{ Row=1135; Fecha=44610; Variedad="Cuatro cascos amarillo"; Calidad="Primera"; Volumen=30;  PMin=22000; PMax=22000; PProm=22000; Unidad="$/caja 15 kilos"; Origen="Región del Maule";             PKg=1467; KgUn=15 },
    @{ Row=1136; Fecha=44610; Variedad="Cuatro cascos verde";    Calidad="Primera"; Volumen=100; PMin=10000; PMax=10000; PProm=10000; Unidad="$/caja 15 kilos"; Origen="Región del Maule";             PKg=667;  KgUn=15 },
    @{ Row=1137; Fecha=44610; Variedad="Morrón rojo";            Calidad="Primera"; Volumen=70;  PMin=20000; PMax=22000; PProm=20857; Unidad="$/caja 18 kilos"; Origen="Provincia de Limarí";          PKg=1159; KgUn=18 },
    @{ Row=1138; Fecha=44610; Variedad="Zafiro rojo";            Calidad="Primera"; Volumen=100; PMin=20000; PMax=20000; PProm=20000; Unidad="$/caja 15 kilos"; Origen="Región de Arica y Parinacota"; PKg=1333; KgUn=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $r.KgUn
    $ws.Cells.Item($row, 18).Value = $clasif
}
